# Update the currency ("Валута") column (C) with the corrected values, fix the
# F4 interest-rate cell that had accidentally been left as placeholder text,
# and tidy up the styling / selection to match the re-uploaded workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> final currency code, derived from the author's corrected data
# (shared strings "ДЕННННН" / "дсадсадса" were replaced by real currency
# codes "ЕУР" (EUR) and "УСД" (USD), and some rows were reassigned to the
# currency that actually matches their interest-rate figures).
$currencyByRow = @{
    2  = "ДЕН"
    3  = "ЕУР"
    4  = "УСД"
    5  = "ДЕН"
    6  = "ЕУР"
    7  = "УСД"
    8  = "ДЕН"
    9  = "ЕУР"
    10 = "УСД"
    11 = "ДЕН"
    12 = "ЕУР"
    13 = "УСД"
    14 = "ДЕН"
    15 = "ЕУР"
    16 = "ДЕН"
    17 = "ЕУР"
    18 = "ДЕН"
    19 = "ЕУР"
    20 = "УСД"
    21 = "ДЕН"
    22 = "ЕУР"
    23 = "УСД"
}

foreach ($row in ($currencyByRow.Keys | Sort-Object)) {
    $ws.Range("C$row").Value = $currencyByRow[$row]
}

# C2 had picked up stray "hyperlink-less" direct formatting (style index 4)
# when the placeholder text was typed in - put it back on the normal body
# style used by every other data cell in the column.
$ws.Range("C2").Style = "Normal 2"

# F4 was accidentally left holding placeholder text ("дсадсадса") instead of
# the nominal interest rate; restore the numeric rate (matching F2/F3) and
# the normal body style.
$ws.Range("F4").Value = 0.01
$ws.Range("F4").Style = "Normal 2"

# Restore the last active selection used when the corrected file was saved.
$ws.Range("C24").Select() | Out-Null
